# "Début de l'ajout des produits à la base de données"
#
# 1) Fix a typo in the existing description on row 15 ("scriptqui" -> "script qui").
# 2) Append a new entry (row 17) to the work-log table "Tableau4" for a work session
#    started on adding products to the database.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Correct the typo in the description of row 15.
$ws.Range("G15").Value = "J'ai fini le script qui crée la base de données"

# 2) Add a new row to the table so it (and the sheet dimension/autofilter) grows
#    from A1:H16 to A1:H17, then fill in its formatting and data.
$lo = $ws.ListObjects.Item(1)
[void]$lo.ListRows.Add()

# Copy the formatting of the previous last row onto the new row.
$ws.Range("A16:H16").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A17").Value = 45061
$ws.Range("B17").Value = 0.33680555555555558
$ws.Range("C17").Value = 0.39583333333333331
$ws.Range("D17").Formula = "=Tableau4[[#This Row],[Heure fin]]-Tableau4[[#This Row],[Heure début]]"
$ws.Range("E17").Value = "CPNV"
$ws.Range("F17").Value = "Base de données"
$ws.Range("G17").Value = "J'ai commencé le script qui ajoutera tous les produits à la base de données"
$ws.Range("H17").Value = "M. Hurni"

# Match the author's final cursor position.
[void]$ws.Range("H17").Select()
